$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.743.81"
$ws.Range("E2").Value = "'  +5.81%  "

$ws.Range("D3").Value = "'2.731.72"
$ws.Range("E3").Value = "'  +4.64%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.00%  "

$ws.Range("D5").Value = "'593.23"
$ws.Range("E5").Value = "'  +1.32%  "

$ws.Range("D6").Value = "'152.46"
$ws.Range("E6").Value = "'  +6.57%  "

$ws.Range("D7").Value = "'0.994"
$ws.Range("E7").Value = "'  -0.44%  "

$ws.Range("D8").Value = "'0.611"
$ws.Range("E8").Value = "'  +2.20%  "

$ws.Range("D9").Value = "'2.763.59"
$ws.Range("E9").Value = "'  +5.48%  "

$ws.Range("D10").Value = "'6.74"
$ws.Range("E10").Value = "'  +3.77%  "

$ws.Range("D11").Value = "'0.113"

$ws.Range("D12").Value = "'0.390"
$ws.Range("E12").Value = "'  +3.41%  "

$ws.Range("E13").Value = "'  +1.64%  "

$ws.Range("D14").Value = "'3.226.29"
$ws.Range("E14").Value = "'  +4.95%  "

$ws.Range("D15").Value = "'26.65"
$ws.Range("E15").Value = "'  +7.25%  "

$ws.Range("D16").Value = "'63.621.86"
$ws.Range("E16").Value = "'  +5.57%  "

$ws.Range("E17").Value = "'  +9.19%  "

$ws.Range("D18").Value = "'2.765.46"
$ws.Range("E18").Value = "'  +5.70%  "

$ws.Range("D19").Value = "'12.06"
$ws.Range("E19").Value = "'  +5.72%  "

$ws.Range("D20").Value = "'4.91"
$ws.Range("E20").Value = "'  +5.98%  "

$ws.Range("D21").Value = "'365.79"
$ws.Range("E21").Value = "'  +5.68%  "

$ws.Range("E22").Value = "'  +1.80%  "

$ws.Range("D23").Value = "'0.538"
$ws.Range("E23").Value = "'  +1.03%  "

$ws.Range("D24").Value = "'0.992"
$ws.Range("E24").Value = "'  -0.59%  "

$ws.Range("D25").Value = "'65.86"
$ws.Range("E25").Value = "'  +3.40%  "

$ws.Range("E26").Value = "'  +5.06%  "

$ws.Range("D27").Value = "'8.70"
$ws.Range("E27").Value = "'  +8.46%  "

$ws.Range("D28").Value = "'0.993"
$ws.Range("E28").Value = "'  -0.50%  "

$ws.Range("D29").Value = "'0.0₃0911"
$ws.Range("E29").Value = "'  +14.49%  "

$ws.Range("E30").Value = "'  +5.54%  "

$ws.Range("D31").Value = "'7.12"
$ws.Range("E31").Value = "'  +9.89%  "

$ws.Range("D32").Value = "'171.72"
$ws.Range("E32").Value = "'  +1.53%  "

$ws.Range("B33").Value = "'Fetch.AI"
$ws.Range("C33").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.20"
$ws.Range("E33").Value = "'  +17.71%  "

$ws.Range("B34").Value = "'USDe"
$ws.Range("C34").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'0.996"
$ws.Range("E34").Value = "'  -0.22%  "

$ws.Range("D35").Value = "'20.67"
$ws.Range("E35").Value = "'  +6.10%  "

$ws.Range("D36").Value = "'4.83"
$ws.Range("E36").Value = "'  +12.49%  "

$ws.Range("E37").Value = "'  +9.48%  "

$ws.Range("D38").Value = "'1.79"
$ws.Range("E38").Value = "'  +9.85%  "

$ws.Range("E39").Value = "'  +19.21%  "

$ws.Range("D40").Value = "'348.87"
$ws.Range("E40").Value = "'  +9.47%  "

$ws.Range("D41").Value = "'4.24"
$ws.Range("E41").Value = "'  +7.96%  "

$ws.Range("D42").Value = "'39.31"
$ws.Range("E42").Value = "'  +2.41%  "

$ws.Range("D43").Value = "'5.67"
$ws.Range("E43").Value = "'  +13.20%  "

$ws.Range("D44").Value = "'22.22"
$ws.Range("E44").Value = "'  +11.38%  "

$ws.Range("B45").Value = "'InjectiveProtocol"
$ws.Range("C45").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'22.28"
$ws.Range("E45").Value = "'  +11.43%  "

$ws.Range("B46").Value = "'Aave"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'143.55"
$ws.Range("E46").Value = "'  +5.85%  "

$ws.Range("B47").Value = "'Hedera"
$ws.Range("C47").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0593"
$ws.Range("E47").Value = "'  +7.88%  "

$ws.Range("B48").Value = "'Mantle"
$ws.Range("C48").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.646"
$ws.Range("E48").Value = "'  +6.23%  "

$ws.Range("D49").Value = "'0.0258"
$ws.Range("E49").Value = "'  +7.09%  "

$ws.Range("E50").Value = "'  +2.56%  "

$ws.Range("D51").Value = "'2.166.31"
$ws.Range("E51").Value = "'  +7.08%  "
